$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Part 1: update the summary sheet "总计" - insert a new first data row
# for 2022-Q3, shifting the existing rows down by one.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Extend column-A styling (bold/border) down to the new row 9 by copying
# the format from the existing last row (A8) before we touch any values.
$summary.Range("A8").Copy()
$summary.Range("A9").PasteSpecial(-4122)

$summaryRows = @(
  @("2022-Q3", 10, 1.24),
  @("2022-Q2", 12, 2.08),
  @("2022-Q1", 11, 1.97),
  @("2021-Q4", 16, 3.04),
  @("2021-Q3", 12, 3.82),
  @("2021-Q2", 8, 3.64),
  @("2021-Q1", 9, 4.38),
  @("2020-Q4", 18, 7.89)
)

for ($i = 0; $i -lt $summaryRows.Count; $i++) {
    $r = $i + 2
    $row = $summaryRows[$i]
    $summary.Cells.Item($r, 1).Value = $i
    $summary.Cells.Item($r, 2).Value = $row[0]
    $summary.Cells.Item($r, 3).Value = $row[1]
    $summary.Cells.Item($r, 4).Value = $row[2]
}

Write-Output "Part 1 (summary sheet) done"

# ---------------------------------------------------------------------
# Part 2: insert the new "2022-Q3" sheet right after "总计" (i.e. right
# before the existing "2022-Q2" tab), with the fund-holdings table.
# ---------------------------------------------------------------------
$beforeTarget = $wb.Worksheets.Item("2022-Q2")
$newSheet = $wb.Worksheets.Add($beforeTarget)
$newSheet.Name = "2022-Q3"

# Reuse the bold/bordered header + index-column formatting from the
# "2022-Q2" sheet so the new sheet matches the existing look. Re-fetch
# the source sheet fresh (by name) right before each Copy - sheet
# references captured before Worksheets.Add()/Name= can go stale.
$fmtSrc = $wb.Worksheets.Item("2022-Q2")
$fmtSrc.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

$fmtSrc = $wb.Worksheets.Item("2022-Q2")
$fmtSrc.Range("A2:A11").Copy()
$newSheet.Range("A2:A11").PasteSpecial(-4122)

Write-Output "Part 2 formatting copied"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $headers.Count; $c++) {
    $newSheet.Cells.Item(1, $c + 2).Value = $headers[$c]
}

$fundRows = @(
  @("001667", "南方转型增长灵活配置混合A", "28.53", "92.86", "1.91", "0.5449", 10),
  @("001044", "嘉实新消费股票", "9.52", "82.84", "5.35", "0.5093", 7),
  @("014499", "南方转型增长灵活配置混合C", "4.30", "92.86", "1.91", "0.0821", 10),
  @("000326", "南方中小盘成长股票", "4.02", "91.10", "1.72", "0.0691", 8),
  @("011179", "浙商智选食品饮料股票A", "0.14", "91.42", "7.17", "0.0100", 4),
  @("011180", "浙商智选食品饮料股票C", "0.08", "91.42", "7.17", "0.0057", 4),
  @("011554", "海富通欣利混合A", "0.79", "34.87", "0.69", "0.0055", 6),
  @("015741", "东财品质生活优选混合A", "0.10", "93.68", "5.15", "0.0052", 10),
  @("011555", "海富通欣利混合C", "0.33", "34.87", "0.69", "0.0023", 6),
  @("015742", "东财品质生活优选混合C", "0.02", "93.68", "5.15", "0.0010", 10)
)

# Force the fund-code and numeric-looking text columns (B, D, E, F, G) to
# stay text, matching the source data (which keeps leading zeros / fixed
# decimal formatting), while A (index) and H (rank) remain real numbers.
$newSheet.Range("B2:B11").NumberFormat = "@"
$newSheet.Range("D2:G11").NumberFormat = "@"

for ($i = 0; $i -lt $fundRows.Count; $i++) {
    $r = $i + 2
    $row = $fundRows[$i]
    $newSheet.Cells.Item($r, 1).Value = $i
    $newSheet.Cells.Item($r, 2).Value = $row[0]
    $newSheet.Cells.Item($r, 3).Value = $row[1]
    $newSheet.Cells.Item($r, 4).Value = $row[2]
    $newSheet.Cells.Item($r, 5).Value = $row[3]
    $newSheet.Cells.Item($r, 6).Value = $row[4]
    $newSheet.Cells.Item($r, 7).Value = $row[5]
    $newSheet.Cells.Item($r, 8).Value = $row[6]
}

Write-Output "done"
